$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.133160666666667
$ws.Range("H2").Value = 24.399482
$ws.Range("I2").Value = 0.3870696756706061
$ws.Range("J2").Value = 0.3870696756706061
$ws.Range("M2").Value = 4.993165333333334
$ws.Range("N2").Value = 14.979496
$ws.Range("O2").Value = 0.06779298131037136
$ws.Range("P2").Value = 0.06779298131037137
$ws.Range("Q2").Value = 40.61021589123023
$ws.Range("R2").Value = 365.491943021072
$ws.Range("S2").Value = 0.0262406072885489
$ws.Range("T2").Value = 0.02624060728854891
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.133160666666667
$ws.Range("H3").Value = 24.399482
$ws.Range("I3").Value = 0.3870696756706061
$ws.Range("J3").Value = 0.3870696756706061
$ws.Range("O3").Value = 0.5355771637189464
$ws.Range("P3").Value = 0.5355771637189464
$ws.Range("Q3").Value = 320.8282601625567
$ws.Range("R3").Value = 2887.45434146301
$ws.Range("S3").Value = 0.2073056790572757
$ws.Range("T3").Value = 0.2073056790572757
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.133160666666667
$ws.Range("H4").Value = 24.399482
$ws.Range("I4").Value = 0.3870696756706061
$ws.Range("J4").Value = 0.3870696756706061
$ws.Range("M4").Value = 29.08216166666666
$ws.Range("N4").Value = 87.24648499999999
$ws.Range("O4").Value = 0.3948530262300277
$ws.Range("P4").Value = 0.3948530262300277
$ws.Range("Q4").Value = 236.5298933689744
$ws.Range("R4").Value = 2128.76904032077
$ws.Range("S4").Value = 0.1528356328004142
$ws.Range("T4").Value = 0.1528356328004142
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 8.133160666666667
$ws.Range("H5").Value = 24.399482
$ws.Range("I5").Value = 0.3870696756706061
$ws.Range("J5").Value = 0.3870696756706061
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.130869
$ws.Range("N5").Value = 0.392607
$ws.Range("O5").Value = 0.001776828740654623
$ws.Range("P5").Value = 0.001776828740654624
$ws.Range("Q5").Value = 1.064378603286
$ws.Range("R5").Value = 9.579407429574001
$ws.Range("S5").Value = 0.0006877565243673966
$ws.Range("T5").Value = 0.0006877565243673967
$ws.Range("I6").Value = 0.3559882250904906
$ws.Range("J6").Value = 0.3559882250904906
$ws.Range("M6").Value = 4.993165333333334
$ws.Range("N6").Value = 14.979496
$ws.Range("O6").Value = 0.06779298131037136
$ws.Range("P6").Value = 0.06779298131037137
$ws.Range("Q6").Value = 37.34924119440267
$ws.Range("R6").Value = 336.143170749624
$ws.Range("S6").Value = 0.0241335030902719
$ws.Range("T6").Value = 0.0241335030902719
$ws.Range("I7").Value = 0.3559882250904906
$ws.Range("J7").Value = 0.3559882250904906
$ws.Range("O7").Value = 0.5355771637189464
$ws.Range("P7").Value = 0.5355771637189464
$ws.Range("S7").Value = 0.1906591639113068
$ws.Range("T7").Value = 0.1906591639113068
$ws.Range("I8").Value = 0.3559882250904906
$ws.Range("J8").Value = 0.3559882250904906
$ws.Range("M8").Value = 29.08216166666666
$ws.Range("N8").Value = 87.24648499999999
$ws.Range("O8").Value = 0.3948530262300277
$ws.Range("P8").Value = 0.3948530262300277
$ws.Range("Q8").Value = 217.5366922644683
$ws.Range("R8").Value = 1957.830230380215
$ws.Range("S8").Value = 0.1405630279792365
$ws.Range("T8").Value = 0.1405630279792365
$ws.Range("I9").Value = 0.3559882250904906
$ws.Range("J9").Value = 0.3559882250904906
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.130869
$ws.Range("N9").Value = 0.392607
$ws.Range("O9").Value = 0.001776828740654623
$ws.Range("P9").Value = 0.001776828740654624
$ws.Range("Q9").Value = 0.9789096734370001
$ws.Range("R9").Value = 8.810187060933
$ws.Range("S9").Value = 0.0006325301096754111
$ws.Range("T9").Value = 0.000632530109675411
$ws.Range("G10").Value = 5.398902333333333
$ws.Range("H10").Value = 16.196707
$ws.Range("I10").Value = 0.2569420992389033
$ws.Range("J10").Value = 0.2569420992389034
$ws.Range("M10").Value = 4.993165333333334
$ws.Range("N10").Value = 14.979496
$ws.Range("O10").Value = 0.06779298131037136
$ws.Range("P10").Value = 0.06779298131037137
$ws.Range("Q10").Value = 26.95761196885244
$ws.Range("R10").Value = 242.618507719672
$ws.Range("S10").Value = 0.01741887093155055
$ws.Range("T10").Value = 0.01741887093155056
$ws.Range("G11").Value = 5.398902333333333
$ws.Range("H11").Value = 16.196707
$ws.Range("I11").Value = 0.2569420992389033
$ws.Range("J11").Value = 0.2569420992389034
$ws.Range("O11").Value = 0.5355771637189464
$ws.Range("P11").Value = 0.5355771637189464
$ws.Range("Q11").Value = 212.9701494143483
$ws.Range("R11").Value = 1916.731344729135
$ws.Range("S11").Value = 0.1376123207503639
$ws.Range("T11").Value = 0.1376123207503639
$ws.Range("G12").Value = 5.398902333333333
$ws.Range("H12").Value = 16.196707
$ws.Range("I12").Value = 0.2569420992389033
$ws.Range("J12").Value = 0.2569420992389034
$ws.Range("M12").Value = 29.08216166666666
$ws.Range("N12").Value = 87.24648499999999
$ws.Range("O12").Value = 0.3948530262300277
$ws.Range("P12").Value = 0.3948530262300277
$ws.Range("Q12").Value = 157.0117504805439
$ws.Range("R12").Value = 1413.105754324895
$ws.Range("S12").Value = 0.1014543654503771
$ws.Range("T12").Value = 0.1014543654503771
$ws.Range("G13").Value = 5.398902333333333
$ws.Range("H13").Value = 16.196707
$ws.Range("I13").Value = 0.2569420992389033
$ws.Range("J13").Value = 0.2569420992389034
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.130869
$ws.Range("N13").Value = 0.392607
$ws.Range("O13").Value = 0.001776828740654623
$ws.Range("P13").Value = 0.001776828740654624
$ws.Range("Q13").Value = 0.706548949461
$ws.Range("R13").Value = 6.358940545149
$ws.Range("S13").Value = 0.0004565421066118158
$ws.Range("T13").Value = 0.000456542106611816
